$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.813315749168396
$ws.Range("B1").Value = 2.236294269561768
$ws.Range("C1").Value = 2.390683889389038
$ws.Range("D1").Value = 3.224337100982666
$ws.Range("E1").Value = 1.691534638404846
